# Update "想去人数" (interest count) figures on the "展览" and "全部类型" sheets
# to reflect the latest refresh of the source data (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 3,5,7,8,10,13,14,15,16,17,19,20,21,22
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1432
$ws1.Range("F5").Value  = 118
$ws1.Range("F7").Value  = 11980
$ws1.Range("F8").Value  = 4449
$ws1.Range("F10").Value = 52
$ws1.Range("F13").Value = 2575
$ws1.Range("F14").Value = 1113
$ws1.Range("F15").Value = 172
$ws1.Range("F16").Value = 56
$ws1.Range("F17").Value = 5186
$ws1.Range("F19").Value = 197
$ws1.Range("F20").Value = 540
$ws1.Range("F21").Value = 11399
$ws1.Range("F22").Value = 11413

# Sheet "全部类型" (all types) - rows shifted by +1 after row 13
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1432
$ws4.Range("F5").Value  = 118
$ws4.Range("F7").Value  = 11980
$ws4.Range("F8").Value  = 4449
$ws4.Range("F10").Value = 52
$ws4.Range("F13").Value = 2575
$ws4.Range("F15").Value = 1113
$ws4.Range("F16").Value = 172
$ws4.Range("F17").Value = 56
$ws4.Range("F18").Value = 5186
$ws4.Range("F20").Value = 197
$ws4.Range("F21").Value = 540
$ws4.Range("F22").Value = 11399
$ws4.Range("F23").Value = 11413
